$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.905.60"
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").Value = "'2.786.83"
$ws.Range("E3").Value = '  -0.93%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'357.73"
$ws.Range("E5").Value = '  +1.54%  '
$ws.Range("D6").Value = "'109.16"
$ws.Range("E6").Value = '  -3.37%  '
$ws.Range("D7").Value = "'0.564"
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("D10").Value = "'39.94"
$ws.Range("E10").Value = '  -3.43%  '
$ws.Range("D11").Value = "'0.0851"
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("E12").Value = '  +1.30%  '
$ws.Range("D13").Value = "'19.48"
$ws.Range("E13").Value = '  -1.93%  '
$ws.Range("D14").Value = "'7.57"
$ws.Range("E14").Value = '  -1.86%  '
$ws.Range("D15").Value = "'3.227.85"
$ws.Range("E15").Value = '  -0.72%  '
$ws.Range("D16").Value = "'2.781.26"
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").Value = "'0.943"
$ws.Range("E17").Value = '  +6.63%  '
$ws.Range("D18").Value = "'51.868.03"
$ws.Range("E18").Value = '  +0.80%  '
$ws.Range("D19").Value = "'7.40"
$ws.Range("E19").Value = '  -0.96%  '
$ws.Range("D20").Value = "'3.12"
$ws.Range("E20").Value = '  -1.99%  '
$ws.Range("D21").Value = "'12.99"
$ws.Range("E21").Value = '  -2.36%  '
$ws.Range("D22").Value = "'0.0₃0979"
$ws.Range("E22").Value = '  -1.23%  '
$ws.Range("D23").Value = "'273.89"
$ws.Range("E23").Value = '  +1.23%  '
$ws.Range("D24").Value = "'70.21"
$ws.Range("E24").Value = '  +0.92%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = "'26.68"
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").Value = "'10.18"
$ws.Range("E28").Value = '  -1.05%  '
$ws.Range("D29").Value = "'0.144"
$ws.Range("E29").Value = '  +4.13%  '
$ws.Range("E30").Value = '  -1.27%  '
$ws.Range("E31").Value = '  +4.10%  '
$ws.Range("D32").Value = "'51.51"
$ws.Range("D33").Value = "'34.27"
$ws.Range("E33").Value = '  +0.91%  '
$ws.Range("E34").Value = '  -1.80%  '
$ws.Range("D35").Value = "'0.0843"
$ws.Range("E35").Value = '  +2.61%  '
$ws.Range("D36").Value = "'5.27"
$ws.Range("E36").Value = '  +1.43%  '
$ws.Range("E37").Value = '  +0.17%  '
$ws.Range("E38").Value = '  +0.94%  '
$ws.Range("D39").Value = "'2.00"
$ws.Range("E39").Value = '  -2.64%  '
$ws.Range("D40").Value = "'17.96"
$ws.Range("E40").Value = '  -0.77%  '
$ws.Range("D41").Value = "'2.54"
$ws.Range("E41").Value = '  +1.00%  '
$ws.Range("E42").Value = '  -1.47%  '
$ws.Range("D43").Value = "'2.26"
$ws.Range("E43").Value = '  -1.44%  '
$ws.Range("D44").Value = "'121.75"
$ws.Range("E44").Value = '  -3.65%  '
$ws.Range("D45").Value = "'22.11"
$ws.Range("E45").Value = '  -6.49%  '
$ws.Range("D46").Value = "'2.071.97"
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").Value = "'3.24"
$ws.Range("E47").Value = '  -2.09%  '
$ws.Range("D48").Value = "'2.19"
$ws.Range("E48").Value = '  -4.33%  '
$ws.Range("D49").Value = "'5.72"
$ws.Range("E49").Value = '  +1.47%  '
$ws.Range("D50").Value = "'0.929"
$ws.Range("E50").Value = '  +0.19%  '
$ws.Range("D51").Value = "'8.93"
$ws.Range("E51").Value = '  +0.51%  '
